# The source data feed added a new weekly record for "Rabanito" (Vega Central
# Mapocho de Santiago) on 2023-07-28 (Excel serial 45135). This record needs
# to be inserted as a new row 372, pushing the existing rows 372:468 down to
# 373:469 (dimension grows from A1:R468 to A1:R469).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 372; this shifts the old rows 372-468
# down to 373-469, preserving all of their existing data untouched.
$ws.Rows(372).Insert()

# Populate the newly inserted row 372 with the new record. Columns A, B, C,
# E-K, N, O, Q, R mirror the template of the surrounding "Rabanito" entries
# for Vega Central Mapocho de Santiago / Metropolitana / Provincia de
# Chacabuco, while D, L, M and P carry the new record's own values.
$ws.Range("A372").Value = 9
$ws.Range("B372").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C372").Value = "Metropolitana"
$ws.Range("D372").Value = 45135
$ws.Range("D372").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E372").Value = 13
$ws.Range("F372").Value = 300000001
$ws.Range("G372").Value = "Rabanito"
$ws.Range("H372").Value = "Sin especificar"
$ws.Range("I372").Value = "Primera"
$ws.Range("J372").Value = 7000
$ws.Range("K372").Value = 3000
$ws.Range("L372").Value = 3500
$ws.Range("M372").Value = 3250
$ws.Range("N372").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O372").Value = "Provincia de Chacabuco"
$ws.Range("P372").Value = 32
$ws.Range("Q372").Value = 100
$ws.Range("R372").Value = "Hortaliza"
